$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("F1").Value = "First day - ENEM 2021"
$ws.Range("G1").Value = "Second day - ENEM 2021"

# Update numeric values (rows 2-28), rounding to 2 decimals per diff
$values = @{
    2  = @{ F = 29.85; G = 35.42 }
    3  = @{ F = 20.21; G = 24.8 }
    4  = @{ F = 46.5;  G = 52.19 }
    5  = @{ F = 26.93; G = 30.68 }
    6  = @{ F = 17.81; G = 21.4 }
    7  = @{ F = 20.12; G = 24.52 }
    8  = @{ F = 26.18; G = 30.15 }
    9  = @{ F = 18.3;  G = 21.77 }
    10 = @{ F = 14.35; G = 16.52 }
    11 = @{ F = 33.16; G = 36.61 }
    12 = @{ F = 11.99; G = 14.63 }
    13 = @{ F = 21.05; G = 23.47 }
    14 = @{ F = 22.49; G = 25.35 }
    15 = @{ F = 15.37; G = 18.51 }
    16 = @{ F = 16.8;  G = 20.11 }
    17 = @{ F = 23.58; G = 27.13 }
    18 = @{ F = 17.88; G = 21.17 }
    19 = @{ F = 16.7;  G = 20.96 }
    20 = @{ F = 14.8;  G = 18.68 }
    21 = @{ F = 15.33; G = 18.74 }
    22 = @{ F = 17.98; G = 21.84 }
    23 = @{ F = 20.08; G = 24.54 }
    24 = @{ F = 17.36; G = 21.44 }
    25 = @{ F = 19.65; G = 24.01 }
    26 = @{ F = 20.17; G = 24.55 }
    27 = @{ F = 41.86; G = 45.03 }
    28 = @{ F = 11.11; G = 14.22 }
}

foreach ($row in $values.Keys) {
    $entry = $values[$row]
    $ws.Cells.Item($row, 6).Value = $entry.F
    $ws.Cells.Item($row, 7).Value = $entry.G
}
